$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the data of rows 136 and 137 (the two matches were reordered) ---
# Columns A,B,C,D,E,G,I are identical between the two rows, only
# F,H,J..V (the match-specific data) need to be swapped.

$row136F = $ws.Range("F136").Value()
$row136H = $ws.Range("H136").Value()
$row136J = $ws.Range("J136").Value()
$row136K = $ws.Range("K136").Value()
$row136L = $ws.Range("L136").Value()
$row136M = $ws.Range("M136").Value()
$row136N = $ws.Range("N136").Value()
$row136O = $ws.Range("O136").Value()
$row136P = $ws.Range("P136").Value()
$row136Q = $ws.Range("Q136").Value()
$row136R = $ws.Range("R136").Value()
$row136S = $ws.Range("S136").Value()
$row136T = $ws.Range("T136").Value()
$row136U = $ws.Range("U136").Value()
$row136V = $ws.Range("V136").Value()

$row137F = $ws.Range("F137").Value()
$row137H = $ws.Range("H137").Value()
$row137J = $ws.Range("J137").Value()
$row137K = $ws.Range("K137").Value()
$row137L = $ws.Range("L137").Value()
$row137M = $ws.Range("M137").Value()
$row137N = $ws.Range("N137").Value()
$row137O = $ws.Range("O137").Value()
$row137P = $ws.Range("P137").Value()
$row137Q = $ws.Range("Q137").Value()
$row137R = $ws.Range("R137").Value()
$row137S = $ws.Range("S137").Value()
$row137T = $ws.Range("T137").Value()
$row137U = $ws.Range("U137").Value()
$row137V = $ws.Range("V137").Value()

$ws.Range("F136").Value = $row137F
$ws.Range("H136").Value = $row137H
$ws.Range("J136").Value = $row137J
$ws.Range("K136").Value = $row137K
$ws.Range("L136").Value = $row137L
$ws.Range("M136").Value = $row137M
$ws.Range("N136").Value = $row137N
$ws.Range("O136").Value = $row137O
$ws.Range("P136").Value = $row137P
$ws.Range("Q136").Value = $row137Q
$ws.Range("R136").Value = $row137R
$ws.Range("S136").Value = $row137S
$ws.Range("T136").Value = $row137T
$ws.Range("U136").Value = $row137U
$ws.Range("V136").Value = $row137V

$ws.Range("F137").Value = $row136F
$ws.Range("H137").Value = $row136H
$ws.Range("J137").Value = $row136J
$ws.Range("K137").Value = $row136K
$ws.Range("L137").Value = $row136L
$ws.Range("M137").Value = $row136M
$ws.Range("N137").Value = $row136N
$ws.Range("O137").Value = $row136O
$ws.Range("P137").Value = $row136P
$ws.Range("Q137").Value = $row136Q
$ws.Range("R137").Value = $row136R
$ws.Range("S137").Value = $row136S
$ws.Range("T137").Value = $row136T
$ws.Range("U137").Value = $row136U
$ws.Range("V137").Value = $row136V

# --- Append the new match as row 167 ---
# Copy the formatting of the last existing row (166) down to the new row
# so the number formats / styles (index + date columns) match.
$ws.Range("A166:V166").Copy($ws.Range("A167:V167"))

$ws.Range("A167").Value = 166
$ws.Range("B167").Value = "spain"
$ws.Range("C167").Value = "laliga2"
$ws.Range("D167").Value = "2023-2024"
$ws.Range("E167").Value = 45247.875
$ws.Range("F167").Value = "Valladolid"
$ws.Range("G167").Value = 1
$ws.Range("H167").Value = "Leganes"
$ws.Range("I167").Value = 1
$ws.Range("J167").Value = 2
$ws.Range("K167").Value = "11/11/2023 18:12"
$ws.Range("L167").Value = 2.19
$ws.Range("M167").Value = "17/11/2023 20:59"
$ws.Range("N167").Value = 3.35
$ws.Range("O167").Value = "11/11/2023 18:12"
$ws.Range("P167").Value = 2.99
$ws.Range("Q167").Value = "17/11/2023 20:56"
$ws.Range("R167").Value = 4.36
$ws.Range("S167").Value = "11/11/2023 18:12"
$ws.Range("T167").Value = 4.22
$ws.Range("U167").Value = "17/11/2023 20:59"
$ws.Range("V167").Value = "https://www.betexplorer.com/football/spain/laliga2/valladolid-leganes/x6gVgkyN/"
